$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comp_type_dmg_algo")
$ws.Activate()

# Insert a new "location" column before "beta" (old column G -> new column H)
$ws.Columns("G").Insert()

# Insert a new "recovery_function" column before "recovery_mean"
# (after the previous insert, "recovery_mean" sits at column N)
$ws.Columns("N").Insert()

# Fill in "recovery_function" (new column N) first, then "location" (new
# column G), so new shared-string entries land in the same order as the
# source edit (recovery_function, Normal, location).
$ws.Cells.Item(1, 14).Value = "recovery_function"
$ws.Range("N2:N25").Value = "Normal"

$ws.Cells.Item(1, 7).Value = "location"
$ws.Range("G2:G25").Value = 0

$ws.Range("B1").Select()
